# Generate Report for Handback
# Adds a new handback record (7f2d55f9-1250-4f40-81b2-ae1e03b433b6.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" sheets/tables.

$wb = $excel.ActiveWorkbook

$fileName      = "7f2d55f9-1250-4f40-81b2-ae1e03b433b6.md"
$pathName      = "e2e\7f2d55f9-1250-4f40-81b2-ae1e03b433b6.md"
$ext           = ".md"
$status        = "Handed back: in sync with en-US"
$srcPath       = "e2e"
$priority      = "ht"
$contentDup    = "True"
$toBeLocalized = "True"
$hasMetadata   = "False"

$zhXlf   = "7f2d55f9-1250-4f40-81b2-ae1e03b433b6.a63845cd2fec44be732facb55dc4505d099f4e0c.zh-cn.xlf"
$deXlf   = "7f2d55f9-1250-4f40-81b2-ae1e03b433b6.a63845cd2fec44be732facb55dc4505d099f4e0c.de-de.xlf"

$zhHandoffDate  = "2017-02-17 08:11:47"
$zhHandbackDate = "2017-02-17 08:12:41"
$deHandoffDate  = "2017-02-17 08:12:04"
$deHandbackDate = "2017-02-17 08:13:04"

$overviewDate = "2017-02-17 08:12:04"

$srcCommit = "6f60dc6da847b4a7a5687a8d29d730c3b0b02278"
$zhCommit  = "d99c171940d59042e14ae003a96b49960445009a"
$deCommit  = "b817f784e7de19685ad9d83127a2bffd1509ca1d"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $overviewDate

$wsOverview.Range("B4").Font.Underline = -4142
$wsOverview.Range("B4").Font.Color = 13527326
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$fileName", $null, $null, $pathName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = $srcPath
$wsZh.Range("E4").Value = $priority
$wsZh.Range("F4").Value = $contentDup
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = $fileName
$wsZh.Range("K4").Value = $zhXlf
$wsZh.Range("L4").Value = $zhHandbackDate
$wsZh.Range("M4").Value = ""
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = $toBeLocalized
$wsZh.Range("P4").Value = ""
$wsZh.Range("Q4").Value = $hasMetadata
$wsZh.Range("R4").Value = ""

$wsZh.Range("A4").Font.Underline = -4142
$wsZh.Range("A4").Font.Color = 13527326
$wsZh.Range("J4").Font.Underline = -4142
$wsZh.Range("J4").Font.Color = 13527326
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/$zhCommit/e2e/$fileName", $null, $null, $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("J4"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/$zhCommit/e2e/$fileName", $null, $null, $fileName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = $srcPath
$wsDe.Range("E4").Value = $priority
$wsDe.Range("F4").Value = $contentDup
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $deHandoffDate
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = $fileName
$wsDe.Range("K4").Value = $deXlf
$wsDe.Range("L4").Value = $deHandbackDate
$wsDe.Range("M4").Value = ""
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = $toBeLocalized
$wsDe.Range("P4").Value = ""
$wsDe.Range("Q4").Value = $hasMetadata
$wsDe.Range("R4").Value = ""

$wsDe.Range("A4").Font.Underline = -4142
$wsDe.Range("A4").Font.Color = 13527326
$wsDe.Range("J4").Font.Underline = -4142
$wsDe.Range("J4").Font.Color = 13527326
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/$deCommit/e2e/$fileName", $null, $null, $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("J4"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/$deCommit/e2e/$fileName", $null, $null, $fileName) | Out-Null

Write-Host "Handback report row added for 7f2d55f9-1250-4f40-81b2-ae1e03b433b6.md"
